# "存款" (deposit) sheet: add bank/deposit_type/currency headers and the
# common metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) that the other
# sheets (land/building/car/insurance/debt) already have, and fix two
# amount cells that were stored as text instead of numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---- Row 1: turn the old "duplicate data" row into a real header row ----
$ws.Cells.Item(1,2).Value = "bank"
$ws.Cells.Item(1,3).Value = "deposit_type"
$ws.Cells.Item(1,4).Value = "currency"
$ws.Cells.Item(1,5).Value = "owner"
$ws.Cells.Item(1,6).Value = "total"

# New header cells G1:M1 - copy formatting from an existing header cell
# (B1) and then set the label.
$headerCols = 7,8,9,10,11,12,13
$headerNames = "property_category","category","date","legislator_name","legislator_id","source_file","index"
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $col = $headerCols[$i]
    $ws.Cells.Item(1,2).Copy()
    $ws.Cells.Item(1,$col).PasteSpecial(-4122)
    $ws.Cells.Item(1,$col).Value = $headerNames[$i]
}

# ---- Rows 2-8: fix two numeric cells that were mis-typed as text ----
$ws.Cells.Item(4,6).Value = 8077460
$ws.Cells.Item(8,6).Value = 4500000

# ---- Rows 2-8: add the new metadata columns G:M ----
# property_category is always "deposit", category is always "normal",
# date/legislator_name/legislator_id/source_file are the same on every
# row (same as every other sheet in this workbook), and index mirrors
# column A (the per-row record id).
$indexValues = 65,66,67,68,69,70,71
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r,2).Copy()

    $ws.Cells.Item($r,7).PasteSpecial(-4122)
    $ws.Cells.Item($r,7).Value = "deposit"

    $ws.Cells.Item($r,8).PasteSpecial(-4122)
    $ws.Cells.Item($r,8).Value = "normal"

    $ws.Cells.Item($r,9).NumberFormat = "@"
    $ws.Cells.Item($r,9).Value = "2013-12-26"
    $ws.Cells.Item($r,2).Copy()
    $ws.Cells.Item($r,9).PasteSpecial(-4122)

    $ws.Cells.Item($r,10).PasteSpecial(-4122)
    $ws.Cells.Item($r,10).Value = "林郁方"

    $ws.Cells.Item($r,11).PasteSpecial(-4122)
    $ws.Cells.Item($r,11).Value = 716

    $ws.Cells.Item($r,12).PasteSpecial(-4122)
    $ws.Cells.Item($r,12).Value = "tmp4c8a1"

    $ws.Cells.Item($r,13).PasteSpecial(-4122)
    $ws.Cells.Item($r,13).Value = $indexValues[$r - 2]
}

$excel.CutCopyMode = 0
